$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextCell "D2" "60.860.31"
Set-TextCell "E2" "  +1.29%  "

# Row 3 - Ethereum
Set-TextCell "D3" "2.598.05"
Set-TextCell "E3" "  +0.58%  "

# Row 4 - TetherUSD
Set-TextCell "E4" "  +0.04%  "

# Row 5 - BNB
Set-TextCell "D5" "520.32"
Set-TextCell "E5" "  +2.94%  "

# Row 6 - Solana
Set-TextCell "D6" "154.41"
Set-TextCell "E6" "  +1.08%  "

# Row 7 - USDC
Set-TextCell "E7" "  +0.00%  "

# Row 8 - XRP
Set-TextCell "D8" "0.590"
Set-TextCell "E8" "  +2.49%  "

# Row 9 - Toncoin
Set-TextCell "D9" "6.69"
Set-TextCell "E9" "  +1.14%  "

# Row 10 - Dogecoin
Set-TextCell "D10" "0.105"
Set-TextCell "E10" "  +2.28%  "

# Row 11 - Cardano
Set-TextCell "D11" "0.347"
Set-TextCell "E11" "  +0.44%  "

# Row 12 - TRON
Set-TextCell "E12" "  +1.63%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextCell "D13" "3.055.27"
Set-TextCell "E13" "  +0.67%  "

# Row 14 - WrappedBTC
Set-TextCell "D14" "60.880.87"
Set-TextCell "E14" "  +1.29%  "

# Row 15 - Avalanche
Set-TextCell "D15" "21.68"
Set-TextCell "E15" "  +1.03%  "

# Row 16 - ShibaInu
Set-TextCell "E16" "  +0.74%  "

# Row 17 - WrappedEther
Set-TextCell "D17" "2.602.10"
Set-TextCell "E17" "  +0.46%  "

# Row 18 - Polkadot
Set-TextCell "D18" "4.73"
Set-TextCell "E18" "  -2.06%  "

# Row 19 - BitcoinCash
Set-TextCell "D19" "352.38"
Set-TextCell "E19" "  +2.03%  "

# Row 20 - Chainlink
Set-TextCell "D20" "10.55"
Set-TextCell "E20" "  +1.60%  "

# Row 21 - Uniswap
Set-TextCell "D21" "6.21"
Set-TextCell "E21" "  +1.15%  "

# Row 22 - Dai
Set-TextCell "E22" "  -0.03%  "

# Row 23 - Litecoin
Set-TextCell "D23" "61.04"
Set-TextCell "E23" "  +1.63%  "

# Row 24 - Polygon
Set-TextCell "D24" "0.426"
Set-TextCell "E24" "  +1.32%  "

# Row 25 & 26 swap: Kaspa <-> WrappedeETH
Set-TextCell "B25" "WrappedeETH"
Set-TextCell "C25" "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextCell "D25" "2.718.53"
Set-TextCell "E25" "  +0.74%  "

Set-TextCell "B26" "Kaspa"
Set-TextCell "C26" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextCell "D26" "0.166"
Set-TextCell "E26" "  +0.28%  "

# Row 27 - Binance-PegBSC-USD
Set-TextCell "D27" "1.00"
Set-TextCell "E27" "  +0.25%  "

# Row 28 - PEPE
Set-TextCell "D28" "0.0₃0846"
Set-TextCell "E28" "  +0.58%  "

# Row 29 - InternetComputer(DFINITY)
Set-TextCell "D29" "7.35"
Set-TextCell "E29" "  +0.22%  "

# Row 30 - USDe
Set-TextCell "E30" "  -0.04%  "

# Row 31 - Aptos
Set-TextCell "D31" "6.32"
Set-TextCell "E31" "  +9.92%  "

# Row 32 - EthereumClassic
Set-TextCell "D32" "19.36"
Set-TextCell "E32" "  +0.30%  "

# Row 33 - PancakeSwap
Set-TextCell "E33" "  +3.07%  "

# Row 34 - Monero
Set-TextCell "D34" "149.05"
Set-TextCell "E34" "  -2.85%  "

# Row 35 - NEARProtocol
Set-TextCell "D35" "4.21"
Set-TextCell "E35" "  +6.04%  "

# Row 36 - SuiNetwork
Set-TextCell "D36" "0.933"
Set-TextCell "E36" "  +9.12%  "

# Row 37 - ImmutableX
Set-TextCell "E37" "  +1.36%  "

# Row 38 - Stacks
Set-TextCell "E38" "  +1.97%  "

# Row 39 - Fetch.AI
Set-TextCell "D39" "0.847"
Set-TextCell "E39" "  +0.23%  "

# Row 40 & 41 swap: Filecoin <-> OKB
Set-TextCell "B40" "OKB"
Set-TextCell "C40" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D40" "36.44"
Set-TextCell "E40" "  +1.73%  "

Set-TextCell "B41" "Filecoin"
Set-TextCell "C41" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextCell "D41" "3.78"
Set-TextCell "E41" "  +0.72%  "

# Row 42 - Bittensor
Set-TextCell "D42" "286.81"
Set-TextCell "E42" "  -2.29%  "

# Row 43 & 44 swap: Mantle <-> Stellar
Set-TextCell "B43" "Stellar"
Set-TextCell "C43" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D43" "0.101"
Set-TextCell "E43" "  +1.41%  "

Set-TextCell "B44" "Mantle"
Set-TextCell "C44" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D44" "0.625"
Set-TextCell "E44" "  +1.73%  "

# Row 45 - Hedera
Set-TextCell "D45" "0.0559"
Set-TextCell "E45" "  +0.23%  "

# Row 46 - FirstDigitalUSD
Set-TextCell "E46" "  +0.05%  "

# Row 47 - EnergySwap
Set-TextCell "D47" "19.56"
Set-TextCell "E47" "  -1.04%  "

# Row 48 & 49 swap: RenderToken <-> VeChain
Set-TextCell "B48" "VeChain"
Set-TextCell "C48" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D48" "0.0237"
Set-TextCell "E48" "  +1.68%  "

Set-TextCell "B49" "RenderToken"
Set-TextCell "C49" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D49" "4.85"
Set-TextCell "E49" "  -0.49%  "

# Row 50 - WhiteBITCoin
Set-TextCell "E50" "  +0.20%  "

# Row 51 - InjectiveProtocol
Set-TextCell "D51" "18.98"
Set-TextCell "E51" "  +7.75%  "
